$wb = $excel.ActiveWorkbook

$wsSlurry  = $wb.Worksheets.Item("Slurry")
$wsClimate = $wb.Worksheets.Item("Climate")

# --- Slurry sheet updates ---
# Dose label changes from 5.7 kg/t to 7.5 kg/t
$wsSlurry.Range("B5").Value = "7.5 kg/t"

# Weather-adjusted formula changes (D5)
$wsSlurry.Range("D5").Formula = "=7.9-1.11"

# Slurry keeps a remembered selection even though it is no longer the active tab
$wsSlurry.Range("D5").Select() | Out-Null

# --- Climate sheet updates ---
# New, more precise weather figures for air.temp / wind.2m / rain.rate
$wsClimate.Range("C2").Value = 4.43101207056639
$wsClimate.Range("D2").Value = 4.05891613991413
$wsClimate.Range("E2").Value = 0.0599629009095261

$wsClimate.Range("C3").Value = 8.23645983645984
$wsClimate.Range("D3").Value = 3.84445591865745
$wsClimate.Range("E3").Value = 0.0552119412831931

$wsClimate.Range("C4").Value = 12.4492495309568
$wsClimate.Range("D4").Value = 3.48391526295633
$wsClimate.Range("E4").Value = 0.0702993488962998

$wsClimate.Range("C5").Value = 16.8762259816193
$wsClimate.Range("D5").Value = 3.15624012423227
$wsClimate.Range("E5").Value = 0.105925308296069

$wsClimate.Range("C6").Value = 14.4977479635841
$wsClimate.Range("D6").Value = 3.32276959833633
$wsClimate.Range("E6").Value = 0.128260170445409

# Apply a higher-precision number format + centered alignment to the updated figures
$wsClimate.Range("C2:E6").NumberFormat = "0.0000"
$wsClimate.Range("C2:E6").HorizontalAlignment = -4108

# Climate becomes the active sheet/tab, with G7 selected
$wsClimate.Activate() | Out-Null
$wsClimate.Range("G7").Select() | Out-Null
